$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Update name column (A) from "pravinN" to "pravin rayN"
$ws2.Range("A1").Value = "pravin ray1"
$ws2.Range("A2").Value = "pravin ray2"
$ws2.Range("A3").Value = "pravin ray3"
$ws2.Range("A4").Value = "pravin ray4"
$ws2.Range("A5").Value = "pravin ray5"

# Update country column (D) from country names to numbers
$ws2.Range("D1").Value = "1"
$ws2.Range("D2").Value = "2"
$ws2.Range("D3").Value = "3"
$ws2.Range("D4").Value = "4"
$ws2.Range("D5").Value = "5"

# Update company size column (G) from ranges to numbers
$ws2.Range("G1").Value = "1"
$ws2.Range("G2").Value = "2"
$ws2.Range("G3").Value = "3"
$ws2.Range("G4").Value = "1"
$ws2.Range("G5").Value = "2"

# Switch the active sheet to Sheet2 and set its selection
$ws2.Activate()
$ws2.Range("E18").Select()
